$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data (header row 4 + values rows 5:29) currently lives in columns
# B, F, J, N, R (with 3 blank columns between each data column). Move it
# into contiguous columns A:E so the labels/ticks line up with the chart.
$ws.Range("A4:A29").Value2 = $ws.Range("B4:B29").Value2
$ws.Range("B4:B29").Value2 = $ws.Range("F4:F29").Value2
$ws.Range("C4:C29").Value2 = $ws.Range("J4:J29").Value2
$ws.Range("D4:D29").Value2 = $ws.Range("N4:N29").Value2
$ws.Range("E4:E29").Value2 = $ws.Range("R4:R29").Value2

# Clear out the now-unused columns F through R.
$ws.Range("F4:R29").Clear()

# Leave the selection where the author left it after the edit.
$ws.Range("G4").Select()
